# Update values produced by re-running the KNN imputation algorithm
# (commit message: "Update Name of Algo")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value  = 7.031000000000001
$ws.Range("E3").Value  = 12.483

$ws.Range("B4").Value  = 6.611999999999999

$ws.Range("C6").Value  = -12.52

$ws.Range("B7").Value  = 6.494
$ws.Range("C7").Value  = -12.918

$ws.Range("B8").Value  = 6.216
$ws.Range("C8").Value  = -11.629
$ws.Range("E8").Value  = 12.576

$ws.Range("E9").Value  = 12.714

$ws.Range("A11").Value = -21.615
$ws.Range("D11").Value = -8.548999999999999

$ws.Range("A12").Value = -21.266
$ws.Range("B12").Value = 6.586

$ws.Range("B14").Value = 6.918000000000001
$ws.Range("D14").Value = -7.764

$ws.Range("A15").Value = -21.192

$ws.Range("C19").Value = -12.545
$ws.Range("D19").Value = -7.928
$ws.Range("E19").Value = 12.806

$ws.Range("C21").Value = -12.588
$ws.Range("D21").Value = -7.528

$ws.Range("B22").Value = 6.686

$ws.Range("C24").Value = -12.255

$ws.Range("C25").Value = -12.69
